$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D are stored as text,
# matching the source data (values like "22.479.46" or "1.002" are
# formatted price strings, not numbers).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.479.46"
$ws.Range("E2").Value = "  +0.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.571.98"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.76"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3717"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.99"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3396"
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07549"
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.30"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.043"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.962"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.571.60"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001123"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.68"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06761"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.300"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.39"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.18"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.483.03"
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("E25").Value = "  -1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.625"
$ws.Range("E26").Value = "  -2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.50"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.061"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.18"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.747.62"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.077"
$ws.Range("E32").Value = "  +7.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.246"
$ws.Range("E33").Value = "  +2.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.011"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.772"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02483"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2304"
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("E39").Value = "  -3.11%  "
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.461"
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.34"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6238"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.815"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5870"
$ws.Range("E47").Value = "  -1.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.14"
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.076"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.218"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07331"
$ws.Range("E51").Value = "  +0.22%  "
